$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 corresponds to "Utilities, fuels, and public services" / "JH+JF".
# Update the values per the corrected JF code for utilities.
$ws.Range("C10").Value  = 9.663478434330043
$ws.Range("I10").Value  = 2.2017940289978917
$ws.Range("K10").Value  = 0.0
$ws.Range("P10").Value  = 4.572825479844843
$ws.Range("R10").Value  = 2.9819572676191917
$ws.Range("S10").Value  = 2.080435302990134
$ws.Range("T10").Value  = 2.4271745201551562
$ws.Range("U10").Value  = 2.066837232407834
$ws.Range("V10").Value  = 0.0
$ws.Range("W10").Value  = 2.1671201072813897
$ws.Range("X10").Value  = 4.137128678467276
$ws.Range("Y10").Value  = 6.2545297559142465
$ws.Range("Z10").Value  = 2.451446265356708
$ws.Range("AA10").Value = 0.0
$ws.Range("AI10").Value = 5.719589512720284
$ws.Range("AL10").Value = 3.058239895395497
$ws.Range("AN10").Value = 2.1983266368262413
$ws.Range("AO10").Value = 2.4271745201551562
$ws.Range("AR10").Value = 4.595920166996558
$ws.Range("AS10").Value = 3.077531884554105
$ws.Range("AT10").Value = 1.668627809630534
$ws.Range("AU10").Value = 2.1393809699081876
$ws.Range("AV10").Value = 0.0
$ws.Range("AY10").Value = 4.495520049625364
$ws.Range("BA10").Value = 8.481848004303611
